$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.212.34"
$ws.Range("D3").Value = "3.781.28"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "3.778.99"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "4.412.83"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "3.796.38"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "68.223.34"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "464.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.697"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000148"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  -3.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "3.735.21"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.299"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.64%  "
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "145.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "391.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "2.783.58"
$ws.Range("E51").Value = "  +4.04%  "
